$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138, shifting existing rows 138.. down by one
# (matches the diff: a new price record is inserted before the old row 138,
# and everything that followed shifts down).
$ws.Rows("138:138").Insert()

# Populate the newly inserted row 138 with the new record's values.
$ws.Cells.Item(138, 1).Value = 8
$ws.Cells.Item(138, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(138, 3).Value = "Coquimbo"
$ws.Cells.Item(138, 4).Value = 44566
$ws.Cells.Item(138, 5).Value = 4
$ws.Cells.Item(138, 6).Value = 100114013
$ws.Cells.Item(138, 7).Value = "Zanahoria"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 800
$ws.Cells.Item(138, 11).Value = 5500
$ws.Cells.Item(138, 12).Value = 6000
$ws.Cells.Item(138, 13).Value = 5750
$ws.Cells.Item(138, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(138, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(138, 16).Value = 288
$ws.Cells.Item(138, 17).Value = 20
$ws.Cells.Item(138, 18).Value = "Hortaliza"

# Match the date-style formatting used by column D elsewhere in the sheet.
$ws.Cells.Item(138, 4).NumberFormat = $ws.Cells.Item(139, 4).NumberFormat
